$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at 1117, pushing the existing 1117:1206 block down to 1120:1209
$ws.Range("A1117:A1119").EntireRow.Insert()

# --- Row 1117 ---
$ws.Range("A1117").Value = 10
$ws.Range("B1117").Value = "Vega Modelo de Temuco"
$ws.Range("C1117").Value = "La Araucanía"
$ws.Range("D1117").Value = 45013
$ws.Range("E1117").Value = 9
$ws.Range("F1117").Value = 100112021
$ws.Range("G1117").Value = "Ají"
$ws.Range("H1117").Value = "Amarillo"
$ws.Range("I1117").Value = "Primera"
$ws.Range("J1117").Value = 20
$ws.Range("K1117").Value = 40000
$ws.Range("L1117").Value = 40000
$ws.Range("M1117").Value = 40000
$ws.Range("N1117").Value = "$/caja 15 kilos"
$ws.Range("O1117").Value = "Región de Arica y Parinacota"
$ws.Range("P1117").Value = 2667
$ws.Range("Q1117").Value = 15
$ws.Range("R1117").Value = "Hortaliza"

# --- Row 1118 ---
$ws.Range("A1118").Value = 10
$ws.Range("B1118").Value = "Vega Modelo de Temuco"
$ws.Range("C1118").Value = "La Araucanía"
$ws.Range("D1118").Value = 45013
$ws.Range("E1118").Value = 9
$ws.Range("F1118").Value = 100112021
$ws.Range("G1118").Value = "Ají"
$ws.Range("H1118").Value = "Americana (o)"
$ws.Range("I1118").Value = "Primera"
$ws.Range("J1118").Value = 30
$ws.Range("K1118").Value = 20000
$ws.Range("L1118").Value = 20000
$ws.Range("M1118").Value = 20000
$ws.Range("N1118").Value = "$/caja 25 kilos"
$ws.Range("O1118").Value = "Provincia de Limarí"
$ws.Range("P1118").Value = 800
$ws.Range("Q1118").Value = 25
$ws.Range("R1118").Value = "Hortaliza"

# --- Row 1119 ---
$ws.Range("A1119").Value = 10
$ws.Range("B1119").Value = "Vega Modelo de Temuco"
$ws.Range("C1119").Value = "La Araucanía"
$ws.Range("D1119").Value = 45013
$ws.Range("E1119").Value = 9
$ws.Range("F1119").Value = 100112021
$ws.Range("G1119").Value = "Ají"
$ws.Range("H1119").Value = "Inferno"
$ws.Range("I1119").Value = "Extra"
$ws.Range("J1119").Value = 20
$ws.Range("K1119").Value = 24000
$ws.Range("L1119").Value = 24000
$ws.Range("M1119").Value = 24000
$ws.Range("N1119").Value = "$/caja 15 kilos"
$ws.Range("O1119").Value = "Región de Arica y Parinacota"
$ws.Range("P1119").Value = 1600
$ws.Range("Q1119").Value = 15
$ws.Range("R1119").Value = "Hortaliza"
